# Shopper_2d / shop_item.xlsx
# Add "sprite_path", "amount" and "max_amount" columns (H, I, J) with per-item
# data used for stand item stacking/limits, matching the new shop CSV schema.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (H1:J1) -------------------------------------------------
# Re-use the same cell formatting the other header cells (A1:G1) already use
# (a "0.00" number format) so the new headers share the existing style.
$ws.Range("H1").Value = "sprite_path"
$ws.Range("H1").NumberFormat = "0.00"

$ws.Range("I1").Value = "amount"
$ws.Range("I1").NumberFormat = "0.00"

$ws.Range("J1").Value = "max_amount"
$ws.Range("J1").NumberFormat = "0.00"

# --- Row 2 : Cracker ---------------------------------------------------------
$ws.Range("H2").Value = "res://Asset/ShopItem/Bread.png"
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 15

# --- Row 3 : Water ------------------------------------------------------------
$ws.Range("H3").Value = "res://Asset/ShopItem/flask.png"
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 30

# --- Row 4 : Sword -------------------------------------------------------------
$ws.Range("H4").Value = "res://Asset/ShopItem/test1.png"
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 4

# --- Row 5 : Shield ------------------------------------------------------------
$ws.Range("H5").Value = "res://Asset/ShopItem/test2.png"
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1

# Resize the new sprite_path column to fit its contents (close to the
# author's saved width) instead of leaving it at the old "bestFit" width.
$ws.Columns.Item(8).ColumnWidth = 7.6

# Put the selection where the author left it after typing in the new data.
$ws.Range("H6").Select()
